$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ScoreM")
$ws2 = $wb.Worksheets.Item("ScoreF")

# Add new data to columns L, M, N for rows 2-11 on ScoreM
$values = @(
    @(35, 31, 10000),
    @(36, 29, 10000),
    @(34, 24, 10000),
    @(34, 57, 10000),
    @(32, 41, 10000),
    @(32, 38, 10000),
    @(33, 33, 10000),
    @(34, 31, 10000),
    @(33, 25, 10000),
    @(35, 13, 10000)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 12).Value = $values[$i][0]
    $ws1.Cells.Item($row, 13).Value = $values[$i][1]
    $ws1.Cells.Item($row, 14).Value = $values[$i][2]
}

# Update selections
$ws1.Range("H7").Select()
$ws2.Range("B25").Select()

# Activate ScoreM (making it the active tab)
$ws1.Activate()
